$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1311.75
$ws.Range("I106").Value = 1311.75
$ws.Range("K106").Value = 1311.75
$ws.Range("M106").Value = -680.75

$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 450
$ws.Range("J111").Value = 600
$ws.Range("K111").Value = 1350
$ws.Range("L111").Value = 1800
$ws.Range("M111").Value = 1717
$ws.Range("N111").Value = -7934

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H137").Value = 1855.8846
$ws.Range("I137").Value = 1644.9474
$ws.Range("J137").Value = 2428.4285
$ws.Range("K137").Value = 4934.8422
$ws.Range("L137").Value = 7285.2855
$ws.Range("M137").Value = -2384.8422
$ws.Range("N137").Value = -12385.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1722.6857
$ws.Range("I32").Value = 1602.7246
$ws.Range("K32").Value = 1602.7246
$ws.Range("M32").Value = -1315.7246

$ws.Range("H45").Value = 1274.1666
$ws.Range("I45").Value = 1274.1666
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1274.1666
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -897.1666
$ws.Range("N45").ClearContents()

$ws.Range("H74").Value = 1515.6666
$ws.Range("I74").Value = 1515.6666
$ws.Range("K74").Value = 1515.6666
$ws.Range("M74").Value = -641.6666

$ws.Range("H77").Value = 1515.6666
$ws.Range("I77").Value = 1515.6666
$ws.Range("K77").Value = 7578.333000000001
$ws.Range("M77").Value = -3210.333000000001

$ws.Range("H110").Value = 757.2
$ws.Range("I110").Value = 395.5
$ws.Range("J110").Value = 998.3333
$ws.Range("K110").Value = 395.5
$ws.Range("L110").Value = 998.3333
$ws.Range("M110").Value = 1649.5
$ws.Range("N110").Value = -5088.3333

$ws.Range("H132").Value = 2499.3125
$ws.Range("I132").Value = 2466
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 7398
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -4868
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4981.091
$ws.Range("I20").Value = 4532.6665
$ws.Range("J20").Value = 6999
$ws.Range("K20").Value = 4532.6665
$ws.Range("L20").Value = 6999
$ws.Range("M20").Value = -4285.6665
$ws.Range("N20").Value = -7493

$ws.Range("H55").Value = 99999
$ws.Range("J55").Value = 99999
$ws.Range("L55").Value = 99999
$ws.Range("N55").Value = -100545

$ws.Range("H86").Value = 3295.1667
$ws.Range("I86").Value = 1880.25
$ws.Range("K86").Value = 1880.25
$ws.Range("M86").Value = -757.25

$ws.Range("H89").Value = 3295.1667
$ws.Range("I89").Value = 1880.25
$ws.Range("K89").Value = 9401.25
$ws.Range("M89").Value = -3785.25

$ws.Range("H94").Value = 1054.6666
$ws.Range("I94").Value = 1088.6875
$ws.Range("K94").Value = 1088.6875
$ws.Range("M94").Value = -637.6875

$ws.Range("H134").Value = 13999.5
$ws.Range("I134").Value = 15000
$ws.Range("J134").Value = 12999
$ws.Range("K134").Value = 45000
$ws.Range("L134").Value = 38997
$ws.Range("M134").Value = -42465
$ws.Range("N134").Value = -44067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1574

$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1566.2646
$ws.Range("I4").Value = 1229.8636
$ws.Range("K4").Value = 3689.5908
$ws.Range("M4").Value = -3577.5908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1672.5
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 2345
$ws.Range("K70").Value = 1000
$ws.Range("L70").Value = 2345
$ws.Range("M70").Value = -730
$ws.Range("N70").Value = -2885

$ws.Range("H73").Value = 1672.5
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 2345
$ws.Range("K73").Value = 1000
$ws.Range("L73").Value = 2345
$ws.Range("M73").Value = -64
$ws.Range("N73").Value = -4217

$ws.Range("H97").Value = 579.375
$ws.Range("J97").Value = 480.25
$ws.Range("L97").Value = 480.25
$ws.Range("N97").Value = -1472.25

$ws.Range("H132").Value = 2787.7273
$ws.Range("I132").Value = 2787.7273
$ws.Range("K132").Value = 8363.1819
$ws.Range("M132").Value = -5833.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19651.611
$ws.Range("I7").Value = 20542.883
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 20542.883
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -20430.883
$ws.Range("N7").Value = -4724

$ws.Range("H55").Value = 200.0625
$ws.Range("I55").Value = 166
$ws.Range("K55").Value = 166
$ws.Range("M55").Value = 7

$ws.Range("H82").Value = 674.5
$ws.Range("I82").Value = 674.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 674.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -313.5
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 674.5
$ws.Range("I85").Value = 674.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 674.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 573.5
$ws.Range("N85").ClearContents()

$ws.Range("H126").Value = 19651.611
$ws.Range("I126").Value = 20542.883
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 61628.649
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -59158.649
$ws.Range("N126").Value = -18440

$ws.Range("H132").Value = 14888.5
$ws.Range("J132").Value = 19777
$ws.Range("L132").Value = 59331
$ws.Range("N132").Value = -64391

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1249.1666
$ws.Range("I126").Value = 498.66666
$ws.Range("J126").Value = 1999.6666
$ws.Range("K126").Value = 1495.99998
$ws.Range("L126").Value = 5998.9998
$ws.Range("M126").Value = 974.0000199999999
$ws.Range("N126").Value = -10938.9998

$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 4999
$ws.Range("K132").Value = 14997
$ws.Range("M132").Value = -12467
